$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
}

Set-TextValue "D2" "303.23"
Set-TextValue "E2" "3.72%"
Set-TextValue "D3" "35.70"
Set-TextValue "E3" "15.10%"
Set-TextValue "D4" "5.089"
Set-TextValue "E4" "2.53%"
Set-TextValue "D5" "0.07815"
Set-TextValue "E5" "4.68%"
Set-TextValue "E6" "1.78%"
Set-TextValue "D7" "8.144"
Set-TextValue "E7" "4.83%"
Set-TextValue "D8" "4.005"
Set-TextValue "E8" "6.36%"
Set-TextValue "D9" "0.9289"
Set-TextValue "E9" "0.87%"
Set-TextValue "D10" "0.09962"
Set-TextValue "E10" "7.23%"
Set-TextValue "D11" "0.1830"
Set-TextValue "E11" "5.57%"
Set-TextValue "D12" "0.08673"
Set-TextValue "E12" "4.44%"
Set-TextValue "D13" "0.03416"
Set-TextValue "E13" "4.45%"
Set-TextValue "D14" "0.09917"
Set-TextValue "E14" "-0.10%"
Set-TextValue "D15" "0.001494"
Set-TextValue "E15" "-0.11%"
Set-TextValue "D16" "0.005600"
Set-TextValue "E16" "-3.26%"
Set-TextValue "D17" "3.482"
Set-TextValue "E17" "0.39%"
Set-TextValue "E18" "-3.94%"
Set-TextValue "D19" "0.3431"
Set-TextValue "E19" "3.00%"
Set-TextValue "D20" "0.1312"
Set-TextValue "E20" "0.96%"
Set-TextValue "D21" "4.580"
Set-TextValue "E21" "10.52%"
Set-TextValue "D22" "0.2233"
Set-TextValue "E22" "5.39%"
Set-TextValue "D23" "0.04665"
Set-TextValue "E23" "3.44%"
Set-TextValue "E24" "1.44%"
Set-TextValue "E25" "5.18%"
Set-TextValue "D26" "0.0001303"
Set-TextValue "E26" "0.50%"
Set-TextValue "D27" "0.0002697"
Set-TextValue "D39" "0.01763"
Set-TextValue "E39" "8.44%"
Set-TextValue "D40" "0.04697"
Set-TextValue "E40" "2.76%"
Set-TextValue "D41" "0.007842"
Set-TextValue "E41" "5.29%"
Set-TextValue "D42" "0.1415"
Set-TextValue "E42" "4.19%"
Set-TextValue "D43" "0.008511"
Set-TextValue "E43" "-13.33%"
Set-TextValue "D44" "0.002216"
Set-TextValue "E44" "2.91%"
Set-TextValue "D45" "0.009214"
Set-TextValue "E45" "-0.48%"
Set-TextValue "D46" "0.00006119"
Set-TextValue "E46" "0.34%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.46%"
Set-TextValue "D48" "5.786"
Set-TextValue "E48" "97.12%"
Set-TextValue "D49" "0.002687"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.46%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.46%"
